# Update column F ("dSF") values for the márquez_germán 2021 save_data sheet.
# This reflects a repull of data where the dSF (change in strikeouts/stat
# differential, final) values are recalculated per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -4
    3  = 2
    4  = -4
    5  = 4
    7  = -4
    8  = -3
    9  = 1
    10 = -7
    11 = 11
    12 = 3
    13 = -3
    14 = -2
    15 = 3
    16 = -1
    17 = 8
    18 = 4
    19 = 5
    20 = -4
    21 = 1
    22 = 1
    23 = -1
    24 = 6
    25 = 4
    26 = -2
    27 = -5
    28 = -4
    30 = 5
    31 = -3
    33 = 1
    34 = -3
    35 = -1
    37 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
